# Apply manifest update:
#  - The previous CommitID ("2475fa8") moves down into a new row 3 (A3),
#    keeping the same cell style used by the CommitID column.
#  - Row 2 gets a new CommitID ("a910999").
#  - DeploymentMethod in row 2 (D2) changes from "StandAlone" to
#    "StandAlone_FirstRun".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the old CommitID value before it gets overwritten.
$oldCommitId = $ws.Range("A2").Value2

# Copy A2's formatting down to the new A3 cell (Style assignment alone does
# not transfer the custom font color used for CommitID cells in this
# runtime, so use copy/paste-special of formats instead).
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = $oldCommitId

# Update row 2 with the new CommitID and the renamed DeploymentMethod value.
$ws.Range("A2").Value = "a910999"
$ws.Range("D2").Value = "StandAlone_FirstRun"

# Update the selected cell to match the target view state.
$ws.Range("D2").Select()
